# edit.ps1 - apply CV content updates per commit diff
# Uses Find to locate a unique text anchor, then assigns Range.Text directly
# (rather than passing replacement text through Find.Execute) so that
# straight quotes/apostrophes are not auto-corrected into curly quotes.
# NOTE: this runtime's PowerShell subset does not bind named (-Param value)
# arguments correctly, so all helper functions below take POSITIONAL args.

$d = $word.ActiveDocument

function Set-RangeTextByFind($FindText, $NewText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Could not find anchor text: $FindText"
    }
    $rng.Text = $NewText
}

# Adds a new bullet paragraph right after the paragraph that currently
# contains $AfterText, cloning that paragraph's formatting (style/numPr),
# and sets its text to $NewBulletText. Returns the new Paragraph object.
# NOTE: Paragraph.Range.Text already includes the trailing paragraph mark
# (chr 13), so concatenating another "`r" onto it and writing it back
# would insert TWO new paragraph marks. Using InsertParagraphAfter()
# directly avoids that double-counting and creates exactly one new
# paragraph, inheriting the source paragraph's formatting.
function Add-BulletAfter($AfterText, $NewBulletText) {
    $paras = $d.Paragraphs
    $count = $paras.Count
    $targetIndex = -1
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$AfterText*") {
            $targetIndex = $i
            break
        }
    }
    if ($targetIndex -eq -1) {
        throw "Could not find paragraph containing: $AfterText"
    }
    $target = $paras.Item($targetIndex)
    $target.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = $NewBulletText
    return $newPara
}

# ---------------------------------------------------------------------
# 1. Professional Summary
# ---------------------------------------------------------------------
Set-RangeTextByFind `
    "Senior Engineering Leader with 15+ years of experience bridging fundamental AI research and enterprise-scale system delivery. Currently leading Google Cloud's Data & Analytics practice for Southeast Asia while driving internal innovations on LLM inference efficiency, multi-agent systems, and AI safety research (activation probing, sandbagging detection).  Proven track record of operating as a ""Player-Coach"": managing regional engineering portfolios while simultaneously architecting and patenting novel frameworks (UPIR, ARTEMIS, FTCS, Speculative Decoding)." `
    'Senior Engineering Leader with 15+ years building the teams, frameworks, and systems that turn Data and AI from research to production. Currently Head of Data & Analytics for Google Cloud in Southeast Asia - a practice built from zero, delivering enterprise Data and AI transformation across 7 countries.  Dual track as "Player-Coach": leading petabyte-scale data platforms and multi-agent systems for Fortune 500 clients, while driving innovation through published research (5 technical disclosures, 6 published packages on PyPI and Maven Central, plus open-source AI safety tools including sandbagging detection and activation steering). Member of Google Cloud delta, architecting solutions at the intersection of applied AI and enterprise scale.'

# ---------------------------------------------------------------------
# 2. Dual-track role paragraph (Google Cloud role summary)
# ---------------------------------------------------------------------
Set-RangeTextByFind `
    'Dual-track role combining technical innovation leadership with regional delivery management. Built Google Cloud''s Data Analytics practice across Southeast Asia while serving as Site Lead overseeing cross-practice operations. Member of <a href="https://cloud.google.com/consulting/innovation-and-transformation">delta</a> - Google Cloud''s innovation and transformation team architecting enterprise AI solutions at scale.' `
    'Dual-track role combining technical innovation leadership with regional delivery management. Built Google Cloud''s Data Analytics practice for Southeast Asia with delivery scope across JAPAC, while serving as Site Lead overseeing cross-practice operations in SEA. Member of <a href="https://cloud.google.com/consulting/innovation-and-transformation">delta</a> - Google Cloud''s innovation and transformation team architecting enterprise AI solutions at scale.'

# ---------------------------------------------------------------------
# 3. Replace $XXM+ placeholders with descriptive language
# ---------------------------------------------------------------------
Set-RangeTextByFind `
    'Direct $XXM+ Data Analytics delivery portfolio across JAPAC while simultaneously overseeing $XXM+ cross-practice portfolio as regional Site Lead.' `
    'Direct regional Data Analytics delivery portfolio across JAPAC while simultaneously overseeing cross-practice portfolio as Site Lead.'

# ---------------------------------------------------------------------
# 4. Remove redundant PSO practice list (already spelled out earlier)
# ---------------------------------------------------------------------
Set-RangeTextByFind `
    'Pioneered agentic AI adoption across all 7 PSO practices (Data Analytics, AI/ML, Infrastructure, Security, Enterprise Architecture, Application Development, Delivery Management) and 6 JAPAC sub-regions, building SDKs, agent catalog, and standardized templates while designing reusable governance frameworks that accelerated innovation and reduced delivery costs.' `
    'Pioneered agentic AI adoption across all 7 PSO practices and 6 JAPAC sub-regions, building SDKs, agent catalog, and standardized templates while designing reusable governance frameworks that accelerated innovation and reduced delivery costs.'

Set-RangeTextByFind `
    'Built agentic tool suites including architecture discovery (100M+ node graph modeling), automated data pipeline generation, and platform cleanup agents that recovered multi-million dollar at-risk engagements and secured significant long-term cloud commitments.' `
    'Built agentic tool suites including architecture discovery (100M+ node graph modeling), automated data pipeline generation, and platform cleanup agents that recovered at-risk engagements and secured significant long-term cloud commitments.'

# ---------------------------------------------------------------------
# 5. Standard Chartered Bank bullets - enhanced platform scale, ML
#    Workbench, MarTech strategy, plus two new bullets (credit risk,
#    data strategy).
# ---------------------------------------------------------------------
Set-RangeTextByFind `
    'Led enterprise-wide AI and data platform development serving 11 markets, delivering technical excellence while influencing C-suite data strategy.' `
    "Led design and development of retail bank's data & analytics platform serving 11 markets, 100+ systems, and 1200+ users."

Set-RangeTextByFind `
    'Delivered a Self-Service ML Platform that reduced model development time from months to weeks' `
    'Developed self-service ML Workbench reducing model deployment time from months to weeks'

Set-RangeTextByFind `
    'MarTech modernization - +30% customer acquisition' `
    'Architected MarTech strategy driving 30% increase in customer acquisition through data-driven personalization'

# Insert two additional Standard Chartered bullets after the MarTech one.
Add-BulletAfter `
    'Architected MarTech strategy driving 30% increase in customer acquisition through data-driven personalization' `
    'Created credit risk models over 15,000+ named entities leveraging news trends and social signals, reducing potential losses by $5M' | Out-Null

Add-BulletAfter `
    'Created credit risk models over 15,000+ named entities leveraging news trends and social signals, reducing potential losses by $5M' `
    'Defined enterprise data strategy including third-party data governance, privacy frameworks, and cloud adoption roadmap' | Out-Null

# ---------------------------------------------------------------------
# 6. Think Big Analytics (Teradata) bullets - data lakes, ad platform,
#    fraud detection, Hadoop clusters.
# ---------------------------------------------------------------------
Set-RangeTextByFind `
    'Data lakes processing 1.2 PB/hour for Fortune 500 clients across APAC' `
    'Designed 5 global data lakes with ETL pipelines handling 1.2 PB/hour and 40K daily files'

Set-RangeTextByFind `
    'Real-time fraud detection systems - 60% reduction in false positives' `
    'Engineered real-time platform processing 2.5M events/second, improving Ad campaign responsiveness by 80%'

Add-BulletAfter `
    'Engineered real-time platform processing 2.5M events/second, improving Ad campaign responsiveness by 80%' `
    'Built ML fraud detection system achieving 60% fewer false positives and 25% higher detection rates, resulting in $3M savings' | Out-Null

Add-BulletAfter `
    'Built ML fraud detection system achieving 60% fewer false positives and 25% higher detection rates, resulting in $3M savings' `
    'Built and managed large-scale Hadoop clusters (300+ nodes) for banks and telcos across JAPAC' | Out-Null

# ---------------------------------------------------------------------
# 7. AI Metacognition Toolkit blurb - mention PyPI publication instead
#    of describing steering vectors inline.
# ---------------------------------------------------------------------
Set-RangeTextByFind `
    'Activation-level detection of sandbagging, deception, and situational awareness in LLMs. Linear probes achieve 90-96% accuracy across Mistral, Gemma, and Qwen models. Includes steering vectors for runtime behavior control.' `
    'Activation-level detection of sandbagging, deception, and situational awareness in LLMs. Linear probes achieve 90-96% accuracy across Mistral, Gemma, and Qwen models. Published on PyPI.'
